# Update "想去人数" (want-to-go count) values in column F across the
# 展览 (exhibition), 演出 (performance) and 全部类型 (all types) sheets.
# These numbers were refreshed from the upstream data source (bilibili
# show platform) at the time the gh-pages site was regenerated.

$wb = $excel.ActiveWorkbook

function Set-FValues($ws, $data) {
    foreach ($row in $data.Keys) {
        $ws.Range("F$row").Value = $data[$row]
    }
}

# 展览 sheet (sheet1)
$wsExhibit = $wb.Worksheets.Item("展览")
Set-FValues $wsExhibit @{
    2  = 15185
    3  = 19569
    5  = 175
    10 = 1
    14 = 216
    15 = 251
    16 = 75
    17 = 1524
    20 = 117
    21 = 249
    22 = 8232
    24 = 46
    25 = 12
    26 = 69
    27 = 1273
    28 = 23
    31 = 6551
    32 = 134
    34 = 189
    36 = 310
    37 = 5602
    38 = 1019
    39 = 31
    41 = 66
}

# 演出 sheet (sheet2)
$wsShow = $wb.Worksheets.Item("演出")
Set-FValues $wsShow @{
    3 = 26
}

# 全部类型 sheet (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
Set-FValues $wsAll @{
    2  = 15185
    3  = 19569
    5  = 175
    10 = 1
    14 = 216
    15 = 251
    16 = 75
    17 = 1524
    21 = 117
    22 = 249
    23 = 8232
    25 = 46
    26 = 12
    27 = 69
    28 = 1273
    29 = 23
    32 = 26
    34 = 6551
    35 = 134
    37 = 189
    39 = 310
    40 = 5602
    41 = 1019
    42 = 31
    44 = 66
}

$wb.Save()
